$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E4").Value = 'RP1-168P16.2;ANKRD11'
$ws.Range("F13").Value = 'protein_coding;antisense'
$ws.Range("E25").Value = 'CTBP1-AS;SPON2'
$ws.Range("F25").Value = 'protein_coding;antisense'
$ws.Range("F38").Value = 'protein_coding;antisense'
$ws.Range("E59").Value = 'RP1-168P16.2;ANKRD11'
$ws.Range("E67").Value = 'RP11-45M22.4;MPRIP;FLCN'
$ws.Range("E76").Value = 'RNF32;LINC01006'
$ws.Range("E77").Value = 'LL22NC03-86G7.1;PPM1F'
$ws.Range("F77").Value = 'protein_coding;antisense'
$ws.Range("F78").Value = 'protein_coding;antisense'
$ws.Range("E80").Value = 'TBCD;ZNF750'
$ws.Range("F83").Value = 'protein_coding;antisense'
$ws.Range("F84").Value = 'protein_coding;antisense'
$ws.Range("E93").Value = 'CTBP1-AS;SPON2'
$ws.Range("F93").Value = 'protein_coding;antisense'
$ws.Range("E94").Value = 'PPP1CA;TBC1D10C'
$ws.Range("E101").Value = 'CTB-147C22.9;KLK6'
$ws.Range("F101").Value = 'protein_coding;antisense'
$ws.Range("F106").Value = 'protein_coding;antisense'
$ws.Range("E119").Value = 'AGAP2;TSPAN31'
$ws.Range("E140").Value = 'TOE1;MUTYH'
$ws.Range("E142").Value = 'ARHGAP9;MARS'
$ws.Range("E146").Value = 'TATDN3;NSL1'
$ws.Range("E150").Value = 'ZNF841;ZNF432'
$ws.Range("F151").Value = 'protein_coding;antisense'
$ws.Range("F153").Value = 'protein_coding;antisense'
$ws.Range("E159").Value = 'RP11-12J10.3;FAM53B'
$ws.Range("E164").Value = 'RP11-33O4.2;NHEJ1'
$ws.Range("E168").Value = 'RP5-892F13.2;TNFRSF9'
$ws.Range("F168").Value = 'protein_coding;processed_pseudogene'
